# Update "想去人数" (F column) values on the 展览 / 演出 / 全部类型 sheets
# to reflect the regenerated gh-pages data output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 245
$ws1.Range("F5").Value = 5846
$ws1.Range("F6").Value = 5173
$ws1.Range("F7").Value = 348
$ws1.Range("F8").Value = 59
$ws1.Range("F12").Value = 32

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 81

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 245
$ws4.Range("F5").Value = 5846
$ws4.Range("F6").Value = 5173
$ws4.Range("F7").Value = 348
$ws4.Range("F8").Value = 59
$ws4.Range("F12").Value = 81
$ws4.Range("F14").Value = 32
